# Updated cryptos list on Wed May  3 09:34:09 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the 50 coin rows (rows 2-51) on the active sheet with new values scraped
# from coinranking.com. Price cells that look like plain numbers are
# written with a leading apostrophe so Excel keeps them as literal text
# (matching how the sheet already stores these values) instead of
# reinterpreting them as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.758.47"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").Value = "1.874.34"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'324.44"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "'0.4620"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").Value = "'0.3869"
$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.07870"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("D10").Value = "'0.9872"
$ws.Range("E10").Value = "  +2.91%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "1.867.07"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").Value = "'6.995"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").Value = "'5.714"
$ws.Range("E14").Value = "  +0.62%  "

$ws.Range("D15").Value = "'0.06992"
$ws.Range("E15").Value = "  +1.83%  "

$ws.Range("D16").Value = "'88.42"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").Value = "'16.81"
$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").Value = "28.775.93"
$ws.Range("E21").Value = "  +2.66%  "

$ws.Range("D22").Value = "'5.284"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "'11.05"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").Value = "'2.103"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "2.127.55"
$ws.Range("E25").Value = "  +4.62%  "

$ws.Range("D26").Value = "'152.64"
$ws.Range("E26").Value = "  -1.40%  "

$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").Value = "'5.833"
$ws.Range("E28").Value = "  +3.28%  "

$ws.Range("D29").Value = "'1.981"
$ws.Range("E29").Value = "  +1.34%  "

$ws.Range("D30").Value = "'119.02"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").Value = "'0.09339"

$ws.Range("D32").Value = "'0.9221"
$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("D33").Value = "'5.308"
$ws.Range("E33").Value = "  +1.16%  "

$ws.Range("D34").Value = "'1.340"
$ws.Range("E34").Value = "  +1.64%  "

$ws.Range("E35").Value = "  +0.43%  "

$ws.Range("D36").Value = "'0.05784"
$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("D37").Value = "'1.154"
$ws.Range("E37").Value = "  +0.68%  "

$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("D39").Value = "'7.678"
$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("D40").Value = "'0.5643"
$ws.Range("E40").Value = "  +1.14%  "

$ws.Range("D41").Value = "'0.1786"
$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("D42").Value = "'9.842"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D43").Value = "'0.07210"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "'11.78"
$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "'2.125"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("E47").Value = "  -2.07%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").Value = "'113.43"
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("E50").Value = "  +3.95%  "

$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.23%  "
